$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "67.940.04"
$c.Style = "Normal"
$ws.Range("E2").Value = "  -1.18%  "

# Row 3
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "3.777.78"
$c.Style = "Normal"
$ws.Range("E3").Value = "  -1.95%  "

# Row 4
$ws.Range("E4").Value = "  +0.11%  "

# Row 5
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "596.58"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.93%  "

# Row 6
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "169.70"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.53%  "

# Row 7
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "3.777.56"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -1.97%  "

# Row 8
$ws.Range("E8").Value = "  +0.06%  "

# Row 9
$ws.Range("E9").Value = "  -0.14%  "

# Row 10
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.164"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -1.77%  "

# Row 11
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "6.55"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +0.65%  "

# Row 12
$ws.Range("E12").Value = "  -1.45%  "

# Row 13
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.0000281"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +0.07%  "

# Row 14
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "36.62"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -0.94%  "

# Row 15
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "4.414.82"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -1.78%  "

# Row 16
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "3.784.94"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -0.77%  "

# Row 17
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "18.88"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +2.93%  "

# Row 18
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "67.906.46"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -1.14%  "

# Row 19
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "7.21"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -2.16%  "

# Row 20
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "0.111"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +0.04%  "

# Row 21
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "10.60"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -3.12%  "

# Row 22
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "467.59"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -1.33%  "

# Row 23
$ws.Range("E23").Value = "  -0.92%  "

# Row 24
$ws.Range("E24").Value = "  -8.09%  "

# Row 25
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "83.72"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +0.16%  "

# Row 26
$ws.Range("E26").Value = "  +0.31%  "

# Row 27
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "10.57"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +2.19%  "

# Row 28
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "12.14"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +0.21%  "

# Row 29
$ws.Range("E29").Value = "  +0.00%  "

# Row 30
$ws.Range("E30").Value = "  -1.24%  "

# Row 31
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "3.929.21"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -1.75%  "

# Row 33
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "30.52"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -2.84%  "

# Row 34
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "2.23"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -2.95%  "

# Row 35
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "9.24"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -0.69%  "

# Row 36
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "3.738.50"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -2.01%  "

# Row 37
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "3.75"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -2.22%  "

# Row 38
$ws.Range("E38").Value = "  +0.09%  "

# Row 39
$ws.Range("E39").Value = "  -0.32%  "

# Row 40
$ws.Range("E40").Value = "  -1.01%  "

# Row 41
$ws.Range("E41").Value = "  -1.36%  "

# Row 42
$ws.Range("E42").Value = "  +0.07%  "

# Row 43
$ws.Range("E43").Value = "  -0.74%  "

# Row 44
$ws.Range("E44").Value = "  -0.01%  "

# Row 45
$ws.Range("B45").Value = "Cosmos"
$ws.Range("C45").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "8.69"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +0.20%  "

# Row 46
$ws.Range("B46").Value = "Stacks"
$ws.Range("C46").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "1.95"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -2.22%  "

# Row 47
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "404.14"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -4.04%  "

# Row 48
$ws.Range("B48").Value = "OKB"
$ws.Range("C48").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "45.63"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -2.68%  "

# Row 49
$ws.Range("B49").Value = "FLOKI"
$ws.Range("C49").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.000281"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -6.62%  "

# Row 50
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "40.18"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +6.47%  "

# Row 51
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "141.30"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -0.22%  "
